# Scheduled market-data refresh: overwrite the crafting-profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for the rows whose
# Universalis price snapshot changed, per job output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 98
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 98
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 98
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -438
# Row 31: Hush Little Wailer / Weak Silencing Potion
$ws.Range("H31").Value = 1357.9
$ws.Range("I31").Value = 215.8
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 647.4000000000001
$ws.Range("L31").Value = 7500
$ws.Range("M31").Value = -417.4000000000001
$ws.Range("N31").Value = -7960
# Row 38: Just Give Him a Serum / Hi-Potion of Strength
$ws.Range("H38").Value = 1040.4348
$ws.Range("I38").Value = 213
$ws.Range("J38").Value = 1676.9231
$ws.Range("K38").Value = 639
$ws.Range("L38").Value = 5030.7693
$ws.Range("M38").Value = -267
$ws.Range("N38").Value = -5774.7693
# Row 39: Riches' Brew / Hi-Potion of Mind
$ws.Range("H39").Value = 228.75
$ws.Range("I39").Value = 56
$ws.Range("J39").Value = 516.6667
$ws.Range("K39").Value = 168
$ws.Range("L39").Value = 1550.0001
$ws.Range("M39").Value = 128
$ws.Range("N39").Value = -2142.0001
# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 3123.3845
$ws.Range("I76").Value = 2668
$ws.Range("J76").Value = 3260
$ws.Range("K76").Value = 2668
$ws.Range("L76").Value = 3260
$ws.Range("M76").Value = -2353
$ws.Range("N76").Value = -3890
# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 3123.3845
$ws.Range("I79").Value = 2668
$ws.Range("J79").Value = 3260
$ws.Range("K79").Value = 2668
$ws.Range("L79").Value = 3260
$ws.Range("M79").Value = -1576
$ws.Range("N79").Value = -5444

$ws = $wb.Worksheets.Item("ARM")
# Row 45: Hollow Hallmarks / Mythril Ingot
$ws.Range("H45").Value = 1134.826
$ws.Range("I45").Value = 843.65
$ws.Range("J45").Value = 3076
$ws.Range("K45").Value = 843.65
$ws.Range("L45").Value = 3076
$ws.Range("M45").Value = -466.65
$ws.Range("N45").Value = -3830
# Row 74: As the Bolt Flies / Titanium Nugget
$ws.Range("H74").Value = 5577.357
$ws.Range("I74").Value = 7000.278
$ws.Range("J74").Value = 3016.1
$ws.Range("K74").Value = 7000.278
$ws.Range("L74").Value = 3016.1
$ws.Range("M74").Value = -6126.278
$ws.Range("N74").Value = -4764.1
# Row 77: Heavy Metal Banned (L) / Titanium Nugget
$ws.Range("H77").Value = 5577.357
$ws.Range("I77").Value = 7000.278
$ws.Range("J77").Value = 3016.1
$ws.Range("K77").Value = 35001.39
$ws.Range("L77").Value = 15080.5
$ws.Range("M77").Value = -30633.39
$ws.Range("N77").Value = -23816.5
# Row 110: Scheduled Maintenance / Deepgold Ingot
$ws.Range("H110").Value = 3885.1155
$ws.Range("I110").Value = 3064.5625
$ws.Range("K110").Value = 3064.5625
$ws.Range("M110").Value = -1019.5625

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 1879.8846
$ws.Range("I134").Value = 1566.5264
$ws.Range("K134").Value = 4699.5792
$ws.Range("M134").Value = -2164.5792

$ws = $wb.Worksheets.Item("CRP")
# Row 94: Beech, Please / Beech Lumber
$ws.Range("H94").Value = 3859.45
$ws.Range("I94").Value = 2211.5
$ws.Range("J94").Value = 4565.7144
$ws.Range("K94").Value = 2211.5
$ws.Range("L94").Value = 4565.7144
$ws.Range("M94").Value = -1760.5
$ws.Range("N94").Value = -5467.7144
# Row 132: Hull Lotta Damage / Ginseng Lumber
$ws.Range("H132").Value = 2649.1052
$ws.Range("I132").Value = 1808.7142
$ws.Range("J132").Value = 5002.2
$ws.Range("K132").Value = 5426.142599999999
$ws.Range("L132").Value = 15006.6
$ws.Range("M132").Value = -2896.142599999999
$ws.Range("N132").Value = -20066.6

$ws = $wb.Worksheets.Item("CUL")
# Row 10: A Real Fungi / Chanterelle Saute
$ws.Range("H10").Value = 223.25
$ws.Range("I10").Value = 223.25
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 669.75
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = -530.75
$ws.Range("N10").ClearContents()
# Row 11: Putting the Squeeze On / Orange Juice
$ws.Range("H11").Value = 157.75
$ws.Range("J11").Value = 300
$ws.Range("L11").Value = 900
$ws.Range("N11").Value = -1180
# Row 12: Butter Me Up / Kukuru Butter
$ws.Range("H12").Value = 14
$ws.Range("I12").Value = 51
$ws.Range("J12").Value = 4.75
$ws.Range("K12").Value = 153
$ws.Range("L12").Value = 14.25
$ws.Range("M12").Value = 20
$ws.Range("N12").Value = -360.25
# Row 13: Fishy Revelations / Braised Pipira
$ws.Range("H13").Value = 4015
$ws.Range("I13").Value = 90
$ws.Range("J13").Value = 4800
$ws.Range("K13").Value = 270
$ws.Range("L13").Value = 14400
$ws.Range("M13").Value = -102
$ws.Range("N13").Value = -14736
# Row 63: The Next to Last Supper / Stuffed Cabbage Rolls
$ws.Range("H63").Value = 2931.1
$ws.Range("I63").Value = 1827.75
$ws.Range("J63").Value = 3666.6667
$ws.Range("K63").Value = 5483.25
$ws.Range("L63").Value = 11000.0001
$ws.Range("M63").Value = -4734.25
$ws.Range("N63").Value = -12498.0001
# Row 66: Nostalgia through the Stomach (L) / Stuffed Cabbage Rolls
$ws.Range("H66").Value = 2931.1
$ws.Range("I66").Value = 1827.75
$ws.Range("J66").Value = 3666.6667
$ws.Range("K66").Value = 16449.75
$ws.Range("L66").Value = 33000.0003
$ws.Range("M66").Value = -12705.75
$ws.Range("N66").Value = -40488.0003
# Row 110: His Dark Utensils / Spaghetti al Nero
$ws.Range("H110").Value = 6800
$ws.Range("J110").Value = 6800
$ws.Range("L110").Value = 20400
$ws.Range("N110").Value = -28580
# Row 114: One Last Meal / Mushroom Saute
$ws.Range("H114").Value = 2321
$ws.Range("I114").Value = 1630.5555
$ws.Range("J114").Value = 3874.5
$ws.Range("K114").Value = 4891.666499999999
$ws.Range("L114").Value = 11623.5
$ws.Range("M114").Value = -1637.666499999999
$ws.Range("N114").Value = -18131.5
# Row 129: Comfort Food / Yakow Moussaka
$ws.Range("H129").Value = 6657.864
$ws.Range("I129").Value = 4752.6665
$ws.Range("J129").Value = 6958.684
$ws.Range("K129").Value = 14257.9995
$ws.Range("L129").Value = 20876.052
$ws.Range("M129").Value = -9257.999500000002
$ws.Range("N129").Value = -30876.052

$ws = $wb.Worksheets.Item("LTW")
# Row 132: Tenets of Tanning / Silver Lobo Leather
$ws.Range("H132").Value = 2540.03
$ws.Range("I132").Value = 1898.5217
$ws.Range("J132").Value = 3967.9033
$ws.Range("K132").Value = 5695.5651
$ws.Range("L132").Value = 11903.7099
$ws.Range("M132").Value = -3165.5651
$ws.Range("N132").Value = -16963.7099
# Row 136: Respect for Br'aax / Br'aax Leather
$ws.Range("H136").Value = 2501.7659
$ws.Range("I136").Value = 1234.3704
$ws.Range("J136").Value = 4212.75
$ws.Range("K136").Value = 3703.1112
$ws.Range("L136").Value = 12638.25
$ws.Range("M136").Value = -1153.1112
$ws.Range("N136").Value = -17738.25

$ws = $wb.Worksheets.Item("WVR")
# Row 107: Flax Wax / Bright Linen Yarn
$ws.Range("H107").Value = 446.25
$ws.Range("I107").Value = 294
$ws.Range("K107").Value = 882
$ws.Range("M107").Value = 1038
# Row 132: Comfy Cabins / Snow Cotton Cloth
$ws.Range("H132").Value = 2393
$ws.Range("I132").Value = 1252
$ws.Range("J132").Value = 3534
$ws.Range("K132").Value = 3756
$ws.Range("L132").Value = 10602
$ws.Range("M132").Value = -1226
$ws.Range("N132").Value = -15662
